# vignette_officer.docx edit:
#  - bold the two "bad" occurrences in the "... automatic testing is bad
#    and I should feel bad." sentence (without touching the unrelated
#    "Automatic testing is bad" table caption elsewhere in the doc)
#  - update the "manual minus auto" confidence interval figures in the
#    auto-generated comparison table

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the target sentence precisely (exact phrase, case sensitive)
#    so we don't accidentally touch the similar "Automatic testing is
#    bad" table caption elsewhere in the document.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "However, automatic testing is bad and I should feel bad."
$find.MatchCase = $true
$find.MatchWholeWord = $false
$found = $find.Execute()

if ($found) {
    $sentenceRange = $d.Range($find.Parent.Start, $find.Parent.End)
    $sentenceEnd = $sentenceRange.End

    # Walk the sentence looking for the standalone word "bad" and bold
    # each match; stay bounded to the sentence so nothing outside it is
    # touched.
    $searchRange = $d.Range($sentenceRange.Start, $sentenceEnd)
    $subFind = $searchRange.Find
    $subFind.ClearFormatting()
    $subFind.Text = "bad"
    $subFind.MatchCase = $true
    $subFind.MatchWholeWord = $true

    while ($searchRange.Start -lt $sentenceEnd -and $subFind.Execute()) {
        if ($searchRange.Start -ge $sentenceEnd -or $searchRange.End -gt $sentenceEnd) {
            break
        }
        $searchRange.Font.Bold = $true
        $searchRange.SetRange($searchRange.End, $sentenceEnd)
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the reported confidence interval for the bootstrap
#    "Displacement" difference-in-means row.
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Text = "manual minus auto: -146.85 [-214.52 to -79.18]"
$find2.Replacement.Text = "manual minus auto: -146.85 [-212.88 to -80.81]"
$find2.MatchCase = $true
$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)

Write-Output "done"
